# Generate Report for Handoff
#
# The handoff process re-ran for the "ac2b4b3f-5c77-4dd6-a812-5e93be2a9ff7.md"
# file, producing a fresh "Latest Handoff Datetime" stamp on both the zh-cn
# and de-de handoff-status sheets (row 4 of each table corresponds to that
# file).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-11-08 22:30:15"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-11-08 22:30:29"
